# Apply the changes described by the diff:
# - Update several odds values in row 2 and row 3
# - Remove row 4 entirely (and the dimension shrinks accordingly)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.4
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 3.2
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("Q2").Value = 2.7
$ws.Range("R2").Value = 1.44
$ws.Range("X2").Value = 10
$ws.Range("Z2").Value = 23
$ws.Range("AE2").Value = 17
$ws.Range("AH2").Value = 15
$ws.Range("AI2").Value = 13
$ws.Range("AJ2").Value = 34
$ws.Range("AL2").Value = 41
$ws.Range("AN2").Value = 4.33
$ws.Range("AO2").Value = 15
$ws.Range("AU2").Value = 9
$ws.Range("AZ2").Value = 67
$ws.Range("BA2").Value = 101

# Row 3 updates
$ws.Range("G3").Value = 2.5
$ws.Range("I3").Value = 2.7
$ws.Range("L3").Value = 3.2
$ws.Range("Q3").Value = 1.67
$ws.Range("R3").Value = 2.15
$ws.Range("S3").Value = 1.3
$ws.Range("T3").Value = 3.4
$ws.Range("W3").Value = 11
$ws.Range("AD3").Value = 7.5
$ws.Range("AE3").Value = 12
$ws.Range("AJ3").Value = 26
$ws.Range("AK3").Value = 19
$ws.Range("AM3").Value = 126
$ws.Range("AS3").Value = 101
$ws.Range("AT3").Value = 3.4
$ws.Range("AX3").Value = 13
$ws.Range("BB3").Value = 101

# Remove the last data row (row 4), which deletes the New Mexico / Phoenix
# Rising match entirely and shrinks the sheet dimension to A1:BD3
$ws.Rows.Item(4).Delete()
